# Auto-generated edit script: update cryptos list values (prices + 1h volume %).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number (e.g. "516.39") must be
# force-formatted as Text *before* the value is written, otherwise Excel's
# COM layer auto-converts the string to a numeric value (dropping trailing
# zeros / introducing float noise), which would not match the source data.

$ws.Range("D2").Value = "56.167.73"
$ws.Range("E2").Value = "  +2.71%  "
$ws.Range("D3").Value = "2.316.44"
$ws.Range("E3").Value = "  +1.61%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "516.39"
$ws.Range("E5").Value = "  +2.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.07"
$ws.Range("E6").Value = "  +3.41%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.533"
$ws.Range("E8").Value = "  +0.92%  "
$ws.Range("D9").Value = "2.337.19"
$ws.Range("E9").Value = "  +1.67%  "
$ws.Range("E10").Value = "  +6.08%  "
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.26"
$ws.Range("E12").Value = "  +7.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.339"
$ws.Range("E13").Value = "  -1.13%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.756.46"
$ws.Range("E14").Value = "  +2.63%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.64"
$ws.Range("E15").Value = "  +1.09%  "
$ws.Range("D16").Value = "56.398.93"
$ws.Range("E16").Value = "  +3.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000133"
$ws.Range("E17").Value = "  +1.80%  "
$ws.Range("D18").Value = "2.349.91"
$ws.Range("E18").Value = "  +2.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.36"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.23"
$ws.Range("E20").Value = "  +2.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "320.15"
$ws.Range("E21").Value = "  +4.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.56"
$ws.Range("E22").Value = "  +1.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.36"
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.67%  "
$ws.Range("E26").Value = "  +5.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.70"
$ws.Range("E27").Value = "  +3.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "170.58"
$ws.Range("E28").Value = "  -0.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.21"
$ws.Range("E29").Value = "  +8.22%  "
$ws.Range("D30").Value = "0.0₃0732"
$ws.Range("E30").Value = "  +4.45%  "
$ws.Range("E31").Value = "  +2.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.19"
$ws.Range("E32").Value = "  +2.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.20"
$ws.Range("E33").Value = "  +1.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.994"
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.933"
$ws.Range("E36").Value = "  +2.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.23"
$ws.Range("E37").Value = "  +3.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.95"
$ws.Range("E38").Value = "  +4.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.51"
$ws.Range("E39").Value = "  +7.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.45"
$ws.Range("E40").Value = "  +2.82%  "
$ws.Range("E41").Value = "  +0.96%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "137.58"
$ws.Range("E42").Value = "  +8.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.55"
$ws.Range("E43").Value = "  +4.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "274.14"
$ws.Range("E44").Value = "  +8.74%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.05"
$ws.Range("E45").Value = "  +0.29%  "
$ws.Range("E46").Value = "  +2.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0501"
$ws.Range("E47").Value = "  +0.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.556"
$ws.Range("E48").Value = "  +1.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0215"
$ws.Range("E49").Value = "  +4.21%  "
$ws.Range("E50").Value = "  +0.96%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.65"
$ws.Range("E51").Value = "  +0.88%  "
